$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7125595211982727
$ws.Range("B1").Value = 1.358111619949341
$ws.Range("C1").Value = 4.197614192962646
$ws.Range("D1").Value = 2.225770950317383
$ws.Range("E1").Value = 0.8065648674964905
